# depozitiMoznosti.xlsx - "Add files via upload" edit
#
# The change:
#  - Clears the deposit-row data (B2:H2) but keeps the existing cell
#    styles/formatting, and also clears G3 (keeping its style).
#  - Row 2 becomes taller (126pt, custom height) - presumably to host a
#    note/link that was pasted into column G/H area.
#  - A new, wide column (G, width 89) is introduced for that note.
#  - The sheet view's selection moves from D12 to H2, scrolled so column B
#    is the first visible column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the data in row 2 (B2:H2), keeping cell formatting/styles ---
$ws.Range("B2:H2").ClearContents()

# --- Clear G3 as well (keeping its style) ---
$ws.Range("G3").ClearContents()

# --- Make row 2 taller to fit the (now empty) wrapped note cell ---
$ws.Rows.Item(2).RowHeight = 126

# --- Widen column G (7) to 89 characters ---
$ws.Columns.Item(7).ColumnWidth = 88.16666666666667

# --- Update the view: scroll so column B is first, select H2 ---
$ws.Range("H2").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
